# Show all missing fields or properties, not just first, in fixture upload.
# This adds a new "field 5" column (with a "fun_fact" example value) to the
# "types" sheet of the test fixture workbook, inserted before the existing
# "field 2 : property 1" column, so the test data exercises a missing value
# in more than one column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("types")

# Insert a new column before the existing "field 2 : property 1" column (H),
# shifting it (and everything after it) one column to the right.
$ws.Columns.Item(8).Insert()

# New header for the inserted column.
$ws.Cells.Item(1, 8).Value = "field 5"

# Only the first data row gets a value for the new column; rows 3 and 4 are
# intentionally left blank so the upload test can verify that all missing
# fields/properties are reported, not just the first one.
$ws.Cells.Item(2, 8).Value = "fun_fact"

# Update the active selection to match the resulting sheet shape.
$ws.Range("J7").Select()
